$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet holds a small table of percentages by year (columns D..M = years
# 2010,2013..2021). Two columns for years 2011 and 2012 need to be inserted
# right after 2010, and a new column for year 2022 needs to be appended at
# the end. Row 8 only contains a stray formatted cell at L8 that must NOT
# move, so we cannot use a whole-column Insert (it would shift L8). Instead
# we shift the existing D..M table data (rows 3-7 only) two columns to the
# right with a block copy, then fill in the newly freed columns E:F and the
# new trailing column P.
# ---------------------------------------------------------------------------

# 1) Shift existing years/data (rows 3-7, cols E..M) right by two columns,
#    landing on G..O. Do the value+format copy, then repeat with a
#    formats-only pass (the engine's "paste all" does not reliably carry
#    per-cell formatting across a multi-row paste, only a formats-only
#    pass does) so every destination cell ends up with the same style as
#    its source cell. Row 8 (L8) is outside the copied range, so it is
#    left completely untouched.
$ws.Range("E3:M7").Copy()
$ws.Range("G3:O7").PasteSpecial(-4163)
$ws.Range("E3:M7").Copy()
$ws.Range("G3:O7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Fill the freed columns E:F with the new 2011/2012 year header + data,
#    using the same base formatting as column D.
$ws.Range("D3:D7").Copy()
$ws.Range("E3:F7").PasteSpecial(-4163)
$ws.Range("D3:D7").Copy()
$ws.Range("E3:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("E4").Value2 = 2011
$ws.Range("F4").Value2 = 2012
$ws.Range("E5").Value2 = 89.6
$ws.Range("F5").Value2 = 87.5
$ws.Range("E6").Value2 = 93.3
$ws.Range("F6").Value2 = 93.9
$ws.Range("E7").Value2 = 92.8
$ws.Range("F7").Value2 = 94.1

# These two inserted data columns carry an explicit 0.0 number format.
$ws.Range("E5:F7").NumberFormat = "0.0"

# 3) Append the new trailing column P (year 2022), copying the formatting
#    from the column immediately to its left (O) then setting values.
$ws.Range("O3:O7").Copy()
$ws.Range("P3:P7").PasteSpecial(-4163)
$ws.Range("O3:O7").Copy()
$ws.Range("P3:P7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("P4").Value2 = 2022
$ws.Range("P5").Value2 = 94.2
$ws.Range("P6").Value2 = 96
$ws.Range("P7").Value2 = 97.5

# P6 uses the 0.0 number format as well.
$ws.Range("P6").NumberFormat = "0.0"

# 4) Restore selection roughly where the editor left off after typing the
#    new column of data.
$ws.Range("Q4").Select()
